$wb = $excel.ActiveWorkbook

# --- Sheet1: narrow column F width ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Columns("F").ColumnWidth = 11.85546875

# --- Sheet2: add new bordered row 27 (C27:F27) and move selection ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("C27").Value = "B"
$ws2.Range("D27").Value = "A"
$ws2.Range("E27").Value = "first"
$ws2.Range("F27").Value = "x[1]"

foreach ($addr in @("C27","D27","E27","F27")) {
    $c = $ws2.Range($addr)
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(7).Weight = 2
    $c.Borders.Item(10).LineStyle = 1
    $c.Borders.Item(10).Weight = 2
}

$ws2.Range("F27").Select()
